$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.0051847438433518765
$ws.Cells.Item(2, 4).Value = -0.17887927139383306
$ws.Cells.Item(2, 5).Value = -0.15855540983906433
$ws.Cells.Item(3, 3).Value = 0.0026342403628657335
$ws.Cells.Item(3, 4).Value = 0.067135267561124493
$ws.Cells.Item(3, 5).Value = 0.077461319928304753
$ws.Cells.Item(4, 3).Value = 0.0076556099333703019
$ws.Cells.Item(4, 4).Value = 0.13763978838000152
$ws.Cells.Item(4, 5).Value = 0.1676492856867271
$ws.Cells.Item(5, 3).Value = 0.0036733817881447243
$ws.Cells.Item(5, 4).Value = -0.087706228082055554
$ws.Cells.Item(5, 5).Value = -0.073306809195380546
$ws.Cells.Item(6, 3).Value = 0.0026857854078715716
$ws.Cells.Item(6, 4).Value = 0.047291368129265467
$ws.Cells.Item(6, 5).Value = 0.05781947311755601
$ws.Cells.Item(7, 3).Value = 0.0059599301564315253
$ws.Cells.Item(7, 4).Value = 0.060010869323229486
$ws.Cells.Item(7, 5).Value = 0.08337340983964231
$ws.Cells.Item(8, 3).Value = 0.0053471701888918874
$ws.Cells.Item(8, 4).Value = -0.2444880145366895
$ws.Cells.Item(8, 5).Value = -0.22352745218061987
$ws.Cells.Item(9, 3).Value = 0.0029868444263396817
$ws.Cells.Item(9, 4).Value = 0.18095115846257134
$ws.Cells.Item(9, 5).Value = 0.1926593960227351
$ws.Cells.Item(10, 3).Value = 0.0097855600900077118
$ws.Cells.Item(10, 4).Value = 0.1356287705023744
$ws.Cells.Item(10, 5).Value = 0.17398753508438955
$ws.Cells.Item(11, 3).Value = 0.0036981139429847899
$ws.Cells.Item(11, 4).Value = -0.12359709309531369
$ws.Cells.Item(11, 5).Value = -0.10910072576220667
$ws.Cells.Item(12, 3).Value = 0.0028071393265849962
$ws.Cells.Item(12, 4).Value = 0.093325913964717191
$ws.Cells.Item(12, 5).Value = 0.10432971846094717
$ws.Cells.Item(13, 3).Value = 0.0062006837720826533
$ws.Cells.Item(13, 4).Value = 0.078805939864667413
$ws.Cells.Item(13, 5).Value = 0.10311221897406572
$ws.Cells.Item(14, 3).Value = 0.0043243717923733341
$ws.Cells.Item(14, 4).Value = -0.23823531740852141
$ws.Cells.Item(14, 5).Value = -0.22128405881698329
$ws.Cells.Item(15, 3).Value = 0.0030690948099396606
$ws.Cells.Item(15, 4).Value = 0.2514717202798179
$ws.Cells.Item(15, 5).Value = 0.26350237404020649
$ws.Cells.Item(16, 3).Value = 0.0097840998351285298
$ws.Cells.Item(16, 4).Value = 0.060187378004015038
$ws.Cells.Item(16, 5).Value = 0.098540418481060696
$ws.Cells.Item(17, 3).Value = 0.0030898194385543423
$ws.Cells.Item(17, 4).Value = -0.14041818268900869
$ws.Cells.Item(17, 5).Value = -0.12830629044749847
$ws.Cells.Item(18, 3).Value = 0.0028130138638081496
$ws.Cells.Item(18, 4).Value = 0.13840069630393365
$ws.Cells.Item(18, 5).Value = 0.1494275286059078
$ws.Cells.Item(19, 3).Value = 0.007182476620697543
$ws.Cells.Item(19, 4).Value = 0.078437810224121191
$ws.Cells.Item(19, 5).Value = 0.10659265376337906
$ws.Cells.Item(20, 3).Value = 0.003799887067565661
$ws.Cells.Item(20, 4).Value = -0.21209690089870434
$ws.Cells.Item(20, 5).Value = -0.19720158860974982
$ws.Cells.Item(21, 3).Value = 0.003533722172674517
$ws.Cells.Item(21, 4).Value = 0.29004413526453465
$ws.Cells.Item(21, 5).Value = 0.3038960983277712
$ws.Cells.Item(22, 3).Value = 0.0076394559588215593
$ws.Cells.Item(22, 4).Value = -0.0011073099869003296
$ws.Cells.Item(22, 5).Value = 0.028838864781199014
$ws.Cells.Item(23, 3).Value = 0.0032074836554781216
$ws.Cells.Item(23, 4).Value = -0.14172726476305811
$ws.Cells.Item(23, 5).Value = -0.12915413640584486
$ws.Cells.Item(24, 3).Value = 0.0030436515990830493
$ws.Cells.Item(24, 4).Value = 0.1598473482379896
$ws.Cells.Item(24, 5).Value = 0.17177826553652345
$ws.Cells.Item(25, 3).Value = 0.0064819346952084888
$ws.Cells.Item(25, 4).Value = 0.062563073054192037
$ws.Cells.Item(25, 5).Value = 0.087971837581093948
$ws.Cells.Item(26, 3).Value = 0.0029291122723972547
$ws.Cells.Item(26, 4).Value = -0.18436044153004363
$ws.Cells.Item(26, 5).Value = -0.17287851029077725
$ws.Cells.Item(27, 3).Value = 0.0036954250000925588
$ws.Cells.Item(27, 4).Value = 0.30237181724697693
$ws.Cells.Item(27, 5).Value = 0.31685764496712865
$ws.Cells.Item(28, 3).Value = 0.0078905737457774239
$ws.Cells.Item(28, 4).Value = -0.051653919564299364
$ws.Cells.Item(28, 5).Value = -0.020723379263354061
$ws.Cells.Item(29, 3).Value = 0.0032715495931157964
$ws.Cells.Item(29, 4).Value = -0.14042774591961016
$ws.Cells.Item(29, 5).Value = -0.12760348323288351
$ws.Cells.Item(30, 3).Value = 0.0032997649455768913
$ws.Cells.Item(30, 4).Value = 0.17778967995605779
$ws.Cells.Item(30, 5).Value = 0.19072454499847577
$ws.Cells.Item(31, 3).Value = 0.0075010098574879802
$ws.Cells.Item(31, 4).Value = 0.042718240192764788
$ws.Cells.Item(31, 5).Value = 0.072121713406367097
$ws.Cells.Item(32, 3).Value = 0.0033613571877060496
$ws.Cells.Item(32, 4).Value = -0.15591773439352435
$ws.Cells.Item(32, 5).Value = -0.1427414309573069
$ws.Cells.Item(33, 3).Value = 0.0038620564454327697
$ws.Cells.Item(33, 4).Value = 0.26581347914141368
$ws.Cells.Item(33, 5).Value = 0.28095249138293915
$ws.Cells.Item(34, 3).Value = 0.0081782655104678387
$ws.Cells.Item(34, 4).Value = -0.088697339357284272
$ws.Cells.Item(34, 5).Value = -0.05663906588905581
$ws.Cells.Item(35, 3).Value = 0.0026388229038821517
$ws.Cells.Item(35, 4).Value = -0.12141972871524116
$ws.Cells.Item(35, 5).Value = -0.11107571370341134
$ws.Cells.Item(36, 3).Value = 0.0033159146440763627
$ws.Cells.Item(36, 4).Value = 0.16771215771611614
$ws.Cells.Item(36, 5).Value = 0.18071032853152452
$ws.Cells.Item(37, 3).Value = 0.0067342937949270093
$ws.Cells.Item(37, 4).Value = 0.022946335862391556
$ws.Cells.Item(37, 5).Value = 0.049344331728774379
$ws.Cells.Item(38, 3).Value = 0.0033969859716475106
$ws.Cells.Item(38, 4).Value = -0.14090120155058575
$ws.Cells.Item(38, 5).Value = -0.12758523557865412
$ws.Cells.Item(39, 3).Value = 0.0042886167333513669
$ws.Cells.Item(39, 4).Value = 0.20126515960028343
$ws.Cells.Item(39, 5).Value = 0.21807626066593397
$ws.Cells.Item(40, 3).Value = 0.0083392265735070768
$ws.Cells.Item(40, 4).Value = -0.097098450667484931
$ws.Cells.Item(40, 5).Value = -0.064409220210877471
$ws.Cells.Item(41, 3).Value = 0.0028996019366023143
$ws.Cells.Item(41, 4).Value = -0.10997595349226479
$ws.Cells.Item(41, 5).Value = -0.098609701548483489
$ws.Cells.Item(42, 3).Value = 0.0034469178619957962
$ws.Cells.Item(42, 4).Value = 0.15204109362027957
$ws.Cells.Item(42, 5).Value = 0.16555278857206068
$ws.Cells.Item(43, 3).Value = 0.0080828040459940253
$ws.Cells.Item(43, 4).Value = -0.013384157281499132
$ws.Cells.Item(43, 5).Value = 0.018299911500243904
$ws.Cells.Item(44, 3).Value = 0.0038306482657553959
$ws.Cells.Item(44, 4).Value = -0.13021337886515058
$ws.Cells.Item(44, 5).Value = -0.11519748466276775
$ws.Cells.Item(45, 3).Value = 0.003784285080706367
$ws.Cells.Item(45, 4).Value = 0.086069307798994091
$ws.Cells.Item(45, 5).Value = 0.10090346130547297
$ws.Cells.Item(46, 3).Value = 0.007890396473709713
$ws.Cells.Item(46, 4).Value = -0.056228374791452256
$ws.Cells.Item(46, 5).Value = -0.025298529385581907
$ws.Cells.Item(47, 3).Value = 0.0033138101986740676
$ws.Cells.Item(47, 4).Value = -0.096059444608924807
$ws.Cells.Item(47, 5).Value = -0.083069523083304245
$ws.Cells.Item(48, 3).Value = 0.0046996977894430544
$ws.Cells.Item(48, 4).Value = 0.1179183760905255
$ws.Cells.Item(48, 5).Value = 0.13634088728426394
$ws.Cells.Item(49, 3).Value = 0.0092971174316332002
$ws.Cells.Item(49, 4).Value = -0.046529575420864833
$ws.Cells.Item(49, 5).Value = -0.010085476751690425
$ws.Cells.Item(50, 3).Value = 0.0037635569967321837
$ws.Cells.Item(50, 4).Value = -0.10975330342741664
$ws.Cells.Item(50, 5).Value = -0.095000402673574166
$ws.Cells.Item(51, 3).Value = 0.0034141793247744476
$ws.Cells.Item(51, 4).Value = -0.011550046712626516
$ws.Cells.Item(51, 5).Value = 0.0018333160949389731
$ws.Cells.Item(52, 3).Value = 0.0086041682097044922
$ws.Cells.Item(52, 4).Value = 0.0056927802665452214
$ws.Cells.Item(52, 5).Value = 0.039420564853666698
$ws.Cells.Item(53, 3).Value = 0.0038356805998852701
$ws.Cells.Item(53, 4).Value = -0.077995659621972399
$ws.Cells.Item(53, 5).Value = -0.062960039896439821
$ws.Cells.Item(54, 3).Value = 0.0053744408982343631
$ws.Cells.Item(54, 4).Value = 0.054845223145966758
$ws.Cells.Item(54, 5).Value = 0.075912683660176061
$ws.Cells.Item(55, 3).Value = 0.011021502910017937
$ws.Cells.Item(55, 4).Value = -0.05589800468056709
$ws.Cells.Item(55, 5).Value = -0.012694426529707207
